# Ontologia_DOCUM.xlsx - "Add files via upload" edit
# Adds two new columns (CategoriaRvt / ClasseIfc) to the "Classes" sheet,
# filled with "null" for every data row, matching the format of the
# neighbouring existing columns, and leaves the workbook focused on the
# "Classes" sheet with the newly added range selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Classes")

# --- formatting -------------------------------------------------------
# Header cells (row 1) should look like the other header cells (e.g. L1:Q1).
$ws.Range("L1").Copy() | Out-Null
$ws.Range("X1:Y1").PasteSpecial(-4122) | Out-Null

# Data cells (rows 2-138) should look like the other plain data cells
# (e.g. Q2, which already carries the same style used across the table).
$ws.Range("Q2").Copy() | Out-Null
$ws.Range("X2:Y138").PasteSpecial(-4122) | Out-Null

# --- values -------------------------------------------------------------
$ws.Range("X1").Value = "CategoriaRvt"
$ws.Range("Y1").Value = "ClasseIfc"
$ws.Range("X2:Y138").Value = "null"

# --- view state -----------------------------------------------------------
# The edit was made/saved while "Classes" was the active sheet, with the new
# columns selected.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 23
$excel.ActiveWindow.ScrollRow = 124
$ws.Range("X2:Y138").Select() | Out-Null
